# Update the "Carrera" column: unify Negocios / Mecatrónica / Manufactura
# into a single "Sistemas" entry, per the commit's shared-strings change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Sistemas"
$ws.Range("D3").Value = "Sistemas"
$ws.Range("D4").Value = "Sistemas"

# Update the selected cell to match the saved view state.
$ws.Range("E8").Select()
